$d = $word.ActiveDocument
$d.Content.Find.Execute("PDF", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Word", 2)
